$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "1047459128"
$ws.Range("D16").Value = "ANDREA CAMILA GONZALEZ ARNEDO"
$ws.Range("E16").Value = "2107"
$ws.Range("F16").Value = 29260

$ws.Range("C17").Value = "1047459128"
$ws.Range("D17").Value = "ANDREA CAMILA GONZALEZ ARNEDO"
$ws.Range("E17").Value = "2106"
$ws.Range("F17").Value = 35112

$ws.Range("C18").Value = "1047459128"
$ws.Range("D18").Value = "ANDREA CAMILA GONZALEZ ARNEDO"
$ws.Range("E18").Value = "2105"
$ws.Range("F18").Value = 35112

$ws.Range("C19").Value = "1047459128"
$ws.Range("D19").Value = "ANDREA CAMILA GONZALEZ ARNEDO"
$ws.Range("E19").Value = "2104"
$ws.Range("F19").Value = 35112

$ws.Range("C20").Value = "1047459128"
$ws.Range("D20").Value = "ANDREA CAMILA GONZALEZ ARNEDO"
$ws.Range("E20").Value = "2103"
$ws.Range("F20").Value = 35112

$ws.Range("C21").Value = "1047459128"
$ws.Range("D21").Value = "ANDREA CAMILA GONZALEZ ARNEDO"
$ws.Range("E21").Value = "2102"
$ws.Range("F21").Value = 35112

$ws.Range("C22").Value = "1047459128"
$ws.Range("D22").Value = "ANDREA CAMILA GONZALEZ ARNEDO"
$ws.Range("E22").Value = "2101"
$ws.Range("F22").Value = 35112

$ws.Range("C23").Value = "1044927883"
$ws.Range("D23").Value = "HERMOGENES ZAPATEIRO PAJARO"
$ws.Range("E23").Value = "2107"
$ws.Range("F23").Value = 29260

$ws.Range("C24").Value = "1044927883"
$ws.Range("D24").Value = "HERMOGENES ZAPATEIRO PAJARO"
$ws.Range("E24").Value = "2106"
$ws.Range("F24").Value = 35112

$ws.Range("C25").Value = "1044927883"
$ws.Range("D25").Value = "HERMOGENES ZAPATEIRO PAJARO"
$ws.Range("E25").Value = "2105"
$ws.Range("F25").Value = 35112

$ws.Range("C26").Value = "1044927883"
$ws.Range("D26").Value = "HERMOGENES ZAPATEIRO PAJARO"
$ws.Range("E26").Value = "2104"
$ws.Range("F26").Value = 35112

$ws.Range("C27").Value = "1044927883"
$ws.Range("D27").Value = "HERMOGENES ZAPATEIRO PAJARO"
$ws.Range("E27").Value = "2103"
$ws.Range("F27").Value = 35112

$ws.Range("C28").Value = "1044927883"
$ws.Range("D28").Value = "HERMOGENES ZAPATEIRO PAJARO"
$ws.Range("E28").Value = "2102"
$ws.Range("F28").Value = 35112

$ws.Range("C29").Value = "1044927883"
$ws.Range("D29").Value = "HERMOGENES ZAPATEIRO PAJARO"
$ws.Range("E29").Value = "2101"
$ws.Range("F29").Value = 35112

$ws.Range("C30").Value = "73572193"
$ws.Range("D30").Value = "MIGUEL DIONISIO LEON DALMAU"
$ws.Range("E30").Value = "2107"
$ws.Range("F30").Value = 29260

$ws.Range("C31").Value = "73572193"
$ws.Range("D31").Value = "MIGUEL DIONISIO LEON DALMAU"
$ws.Range("E31").Value = "2106"
$ws.Range("F31").Value = 35112

$ws.Range("C32").Value = "73572193"
$ws.Range("D32").Value = "MIGUEL DIONISIO LEON DALMAU"
$ws.Range("E32").Value = "2105"
$ws.Range("F32").Value = 35112

$ws.Range("C33").Value = "73572193"
$ws.Range("D33").Value = "MIGUEL DIONISIO LEON DALMAU"
$ws.Range("E33").Value = "2104"
$ws.Range("F33").Value = 35112

$ws.Range("C34").Value = "73572193"
$ws.Range("D34").Value = "MIGUEL DIONISIO LEON DALMAU"
$ws.Range("E34").Value = "2103"
$ws.Range("F34").Value = 35112

$ws.Range("C35").Value = "73572193"
$ws.Range("D35").Value = "MIGUEL DIONISIO LEON DALMAU"
$ws.Range("E35").Value = "2102"
$ws.Range("F35").Value = 35112

$ws.Range("C36").Value = "73572193"
$ws.Range("D36").Value = "MIGUEL DIONISIO LEON DALMAU"
$ws.Range("E36").Value = "2101"
$ws.Range("F36").Value = 35112

$ws.Range("C37").Value = "73572193"
$ws.Range("D37").Value = "MIGUEL DIONISIO LEON DALMAU"
$ws.Range("E37").Value = "2012"
$ws.Range("F37").Value = 35112

$ws.Range("C38").Value = "73572193"
$ws.Range("D38").Value = "MIGUEL DIONISIO LEON DALMAU"
$ws.Range("E38").Value = "2011"
$ws.Range("F38").Value = 35112

$ws.Range("C39").Value = "73572193"
$ws.Range("D39").Value = "MIGUEL DIONISIO LEON DALMAU"
$ws.Range("E39").Value = "2010"
$ws.Range("F39").Value = 35112

$ws.Range("C40").Value = "73572193"
$ws.Range("D40").Value = "MIGUEL DIONISIO LEON DALMAU"
$ws.Range("E40").Value = "2009"
$ws.Range("F40").Value = 35112

$ws.Range("C41").Value = "73572193"
$ws.Range("D41").Value = "MIGUEL DIONISIO LEON DALMAU"
$ws.Range("E41").Value = "2008"
$ws.Range("F41").Value = 35112
